# Add four new monthly rows (303-306) to the CNCBBS balance-sheet data,
# continuing the series stored in columns A:G of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 303; Date = 45107; Open = 41806284000000; High = 41806284000000; Low = 41806284000000; Close = 41806284000000; Volume = 0 },
    @{ Row = 304; Date = 45138; Open = 40809168000000; High = 40809168000000; Low = 40809168000000; Close = 40809168000000; Volume = 0 },
    @{ Row = 305; Date = 45169; Open = 41684045000000; High = 41684045000000; Low = 41684045000000; Close = 41684045000000; Volume = 0 },
    @{ Row = 306; Date = 45199; Open = 42735489000000; High = 42735489000000; Low = 42735489000000; Close = 42735489000000; Volume = 0 }
)

$lastExistingRow = 302

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A carries the date style (numFmt "YYYY-MM-DD HH:MM:SS", bold,
    # bordered, centered) used by every data row above it. Copy that
    # formatting down from the previous row, then overwrite the value.
    $srcCell = $ws.Cells.Item($lastExistingRow, 1)
    $dstCell = $ws.Cells.Item($row, 1)
    $srcCell.Copy($dstCell)
    $dstCell.Value = $r.Date

    $ws.Cells.Item($row, 2).Value = "ECONOMICS:CNCBBS"
    $ws.Cells.Item($row, 3).Value = $r.Open
    $ws.Cells.Item($row, 4).Value = $r.High
    $ws.Cells.Item($row, 5).Value = $r.Low
    $ws.Cells.Item($row, 6).Value = $r.Close
    $ws.Cells.Item($row, 7).Value = $r.Volume
}
